$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old, wider report layout (columns D:Q held the
# civilStatus/spouse/employment/etc. fields that this report no longer
# needs) and drop the old per-person duplicate detail rows.
$ws.Cells.Clear()
$ws.Columns("D:Q").Delete()

# New, smaller "weekly/monthly filtered" report: Fullname / Amount / ContactNo
$ws.Range("A1").Value = "Fullname"
$ws.Range("B1").Value = "Amount"
$ws.Range("C1").Value = "ContactNo"

$ws.Range("A2").Value = "Karl Borromeo"
$ws.Range("B2").Value = 1000
$ws.Range("C2").Value = "09123456789"

$ws.Range("A3").Value = "Mary Grace Galllardo"
$ws.Range("B3").Value = 2000
$ws.Range("C3").Value = "09123456789"

$ws.Range("A4").Value = "Kryzz Andig"
$ws.Range("B4").Value = 4000
$ws.Range("C4").Value = "09123456789"

$ws.Range("A5").Value = "Total:"
$ws.Range("B5").Value = 7000

$ws.Columns("A:C").ColumnWidth = 24.166666666666668
